$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.597.78"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "3.583.54"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "609.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").Value = "3.582.15"
$ws.Range("E7").Value = "  +0.52%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "4.191.43"
$ws.Range("E13").Value = "  +0.66%  "
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.03"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "3.579.09"
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").Value = "66.630.01"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("E19").Value = "  -2.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.23%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.622"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("D25").Value = "3.726.95"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.46%  "
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.15%  "
$ws.Range("D32").Value = "3.579.82"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "
$ws.Range("E34").Value = "  -3.09%  "
$ws.Range("E35").Value = "  -2.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.61"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "173.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.31%  "
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("E42").Value = "  -1.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.894"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.68%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.49%  "
$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "24.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.32%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.21"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.80%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.16%  "
